$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.229.80"
$ws.Range("D2").ClearFormats()

$ws.Range("D3").Value = "'1.854.57"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'0.7005"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.25%  "

$ws.Range("D6").Value = "'237.75"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "'0.08080"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.02%  "

$ws.Range("D9").Value = "'0.3018"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.38%  "

$ws.Range("D10").Value = "'23.45"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.72%  "

$ws.Range("D11").Value = "'0.08176"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.29%  "

$ws.Range("D12").Value = "'1.856.82"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.28%  "

$ws.Range("D13").Value = "'5.193"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.17%  "

$ws.Range("D14").Value = "'0.7053"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.20%  "

$ws.Range("D15").Value = "'89.80"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.73%  "

$ws.Range("D16").Value = "'29.262.29"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").Value = "'5.828"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.97%  "

$ws.Range("D18").Value = "'0.000007895"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.42%  "

$ws.Range("D19").Value = "'13.26"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.06%  "

$ws.Range("D20").Value = "'236.98"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.61%  "

$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").Value = "'2.116.78"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.00%  "

$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("D24").Value = "'7.457"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("D25").Value = "'162.89"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.75%  "

$ws.Range("D26").Value = "'8.878"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.70%  "

$ws.Range("E27").Value = "  -0.72%  "

$ws.Range("D28").Value = "'18.05"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.18%  "

$ws.Range("D29").Value = "'1.918"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.88%  "

$ws.Range("E30").Value = "  +0.88%  "

$ws.Range("D31").Value = "'1.474"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.92%  "

$ws.Range("D32").Value = "'4.356"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.37%  "

$ws.Range("D33").Value = "'4.023"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.64%  "

$ws.Range("D34").Value = "'0.05186"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.17%  "

$ws.Range("D35").Value = "'1.161"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.30%  "

$ws.Range("D36").Value = "'0.7185"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.43%  "

$ws.Range("D37").Value = "'0.9990"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.59%  "

$ws.Range("D38").Value = "'2.687"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.81%  "

$ws.Range("D39").Value = "'0.01849"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.25%  "

$ws.Range("D40").Value = "'2.723"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.64%  "

$ws.Range("D41").Value = "'0.9319"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.17%  "

$ws.Range("D42").Value = "'1.145.72"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.18%  "

$ws.Range("D43").Value = "'6.018"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.59%  "

$ws.Range("D44").Value = "'0.4250"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("D45").Value = "'70.22"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.13%  "

$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").Value = "'102.94"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("E48").Value = "  -3.66%  "

$ws.Range("D49").Value = "'1.744"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.81%  "

$ws.Range("D50").Value = "'2.006.40"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("D51").Value = "'9.150"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.16%  "
